$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.313.73'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.871.61'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4669'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2843'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06558'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.20'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07883'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.89'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").Value = '1.866.22'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.114'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6756'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '281.05'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '30.300.55'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.502'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("D21").Value = '2.118.44'
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007288'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.171'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.236'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.929'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.371'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09721'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.422'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.476'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.110'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.121'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7059'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.713'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01860'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.291'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.544'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.23'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.952'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8464'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4178'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.189'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.153'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '932.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1127'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.10%  '
